# Update data: 2025-10-29 10:21
# Applies the refreshed market-health snapshot to the workbook:
#   - Metadata!A2 timestamp bump
#   - Top Losers: rows re-ranked/re-valued as the underlying data refreshed
#   - 1 Month Performance: two rows re-ranked/re-valued

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Last Updated" timestamp ---------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 10:21 AM"

# --- Top Losers sheet: rows 22-23, 27-29, 40-45 refreshed ------------------
$losers = $wb.Worksheets.Item("Top Losers")

$losersData = @{
    22 = @("ATHERENERG", -4.0945, -0.0142, 24.8806)
    23 = @("SHAREINDIA", -4.0806, -1.6889, 54.7217)
    27 = @("360ONE",     -3.8488, -4.976,  10.0293)
    28 = @("SUMMITSEC",  -3.8325, -2.0765, 5.5474)
    29 = @("SMSPHARMA",  -3.7339, -3.0871, 17.4387)
    40 = @("BOSCHLTD",   -3.0099, -3.123,  -2.0055)
    41 = @("DRREDDY",    -2.9859, -2.5475, 2.2228)
    42 = @("ROSSTECH",   -2.9778, 1.9028,  -6.8057)
    43 = @("OAL",        -2.9496, -1.278,  8.7362)
    44 = @("ENDURANCE",  -2.939,  -2.2945, 3.4531)
    45 = @("POLICYBZR",  -2.907,  2.2365,  1.2573)
}

foreach ($row in $losersData.Keys) {
    $vals = $losersData[$row]
    $losers.Cells.Item($row, 2).Value = $vals[0]
    $losers.Cells.Item($row, 3).Value = $vals[1]
    $losers.Cells.Item($row, 4).Value = $vals[2]
    $losers.Cells.Item($row, 5).Value = $vals[3]
}

# --- 1 Month Performance sheet: rows 22-23 refreshed ------------------------
$perf = $wb.Worksheets.Item("1 Month Performance")

$perfData = @{
    22 = @("SHAREINDIA", 35.3207)
    23 = @("SOUTHBANK",  35.2819)
}

foreach ($row in $perfData.Keys) {
    $vals = $perfData[$row]
    $perf.Cells.Item($row, 2).Value = $vals[0]
    $perf.Cells.Item($row, 3).Value = $vals[1]
}
